# edit.ps1 -- applies the "New crime data collected" update to the weekly
# 6th Precinct CompStat sheet (issue number, report week, and refreshed stats).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: issue number 25 -> 26, report week 6/19-6/25 -> 6/26-7/2 ---
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# --- Cells whose type/number-format changes (blank "***.* " <-> numeric) ---
# Write the new value first, then paste-special the number format only from a
# donor cell that already carries the target style, so the stored style id matches.
$ws.Range("C14").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C22").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("D27").Value = 3
$ws.Range("F14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = 33.333333333333
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("D30").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

# --- Refreshed weekly / 28-day / YTD / 2-year figures (value-only updates) ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -50
$ws.Range("M15").Value = -33.333333333333
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -37.5
$ws.Range("F16").Value = 18
$ws.Range("H16").Value = -21.739130434782
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 122
$ws.Range("K16").Value = -22.950819672131
$ws.Range("L16").Value = 8.045977011494
$ws.Range("M16").Value = 42.424242424242
$ws.Range("N16").Value = -77.830188679245
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 87
$ws.Range("K17").Value = 26.436781609195
$ws.Range("L17").Value = 18.279569892473
$ws.Range("M17").Value = 134.042553191489
$ws.Range("N17").Value = -19.117647058823
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = -64
$ws.Range("I18").Value = 153
$ws.Range("J18").Value = 219
$ws.Range("K18").Value = -30.136986301369
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = 54.545454545454
$ws.Range("N18").Value = -61.265822784810
$ws.Range("C19").Value = 37
$ws.Range("D19").Value = 45
$ws.Range("E19").Value = -17.777777777777
$ws.Range("F19").Value = 120
$ws.Range("G19").Value = 134
$ws.Range("H19").Value = -10.447761194029
$ws.Range("I19").Value = 618
$ws.Range("J19").Value = 611
$ws.Range("K19").Value = 1.145662847790
$ws.Range("L19").Value = 103.960396039604
$ws.Range("M19").Value = 19.075144508670
$ws.Range("N19").Value = -47.848101265822
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -37.5
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -41.379310344827
$ws.Range("L20").Value = 6.25
$ws.Range("M20").Value = -15
$ws.Range("N20").Value = -95.442359249329
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 67
$ws.Range("E21").Value = -14.925373134328
$ws.Range("F21").Value = 198
$ws.Range("G21").Value = 242
$ws.Range("H21").Value = -18.181818181818
$ws.Range("I21").Value = 997
$ws.Range("J21").Value = 1077
$ws.Range("K21").Value = -7.428040854224
$ws.Range("L21").Value = 64.521452145214
$ws.Range("M21").Value = 31.704095112285
$ws.Range("N21").Value = -60.420801905518
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("C24").Value = 67
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = 48.888888888888
$ws.Range("F24").Value = 250
$ws.Range("G24").Value = 212
$ws.Range("H24").Value = 17.924528301886
$ws.Range("I24").Value = 1010
$ws.Range("J24").Value = 965
$ws.Range("K24").Value = 4.663212435233
$ws.Range("L24").Value = 68.896321070234
$ws.Range("M24").Value = 27.364438839848
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 40
$ws.Range("H25").Value = -14.893617021276
$ws.Range("I25").Value = 230
$ws.Range("J25").Value = 197
$ws.Range("K25").Value = 16.751269035533
$ws.Range("L25").Value = 76.923076923076
$ws.Range("M25").Value = 82.539682539682
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -50
$ws.Range("L26").Value = -14.285714285714
$ws.Range("C27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 36
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = 9.090909090909
$ws.Range("L27").Value = 44

$excel.CutCopyMode = 0
